# Natmi following Dr Hou advice
# Update A2m-Lrp1 ligand-receptor rows: "ECs" now sends to all target clusters (rows 2-6,
# previously mis-offset so the sending cluster incorrectly showed as "FAPs"/wrong columns),
# and a new sending cluster "FAPs" -> A2m/Lrp1 block is added for all target clusters (rows 7-11).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> A2m/Lrp1 -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "A2m"
$ws.Range("C2").Value = "Lrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1.0
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.007329
$ws.Range("H2").Value = 0.021987
$ws.Range("I2").Value = 0.08410957541630165
$ws.Range("J2").Value = 0.08410957541630165
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 6.744736333333333
$ws.Range("N2").Value = 20.234209
$ws.Range("O2").Value = 0.01049273071342527
$ws.Range("P2").Value = 0.01049273071342527
$ws.Range("Q2").Value = 0.04943217258699999
$ws.Range("R2").Value = 0.444889553283
$ws.Range("S2").Value = 0.0008825391252637876
$ws.Range("T2").Value = 0.0008825391252637877

# Row 3: ECs -> A2m/Lrp1 -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "A2m"
$ws.Range("C3").Value = "Lrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1.0
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.007329
$ws.Range("H3").Value = 0.021987
$ws.Range("I3").Value = 0.08410957541630165
$ws.Range("J3").Value = 0.08410957541630165
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 182.5316823333334
$ws.Range("N3").Value = 547.595047
$ws.Range("O3").Value = 0.283963033503136
$ws.Range("P3").Value = 0.2839630335031361
$ws.Range("Q3").Value = 1.337774699821
$ws.Range("R3").Value = 12.039972298389
$ws.Range("S3").Value = 0.02388401018187381
$ws.Range("T3").Value = 0.02388401018187382

# Row 4: ECs -> A2m/Lrp1 -> M1
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "A2m"
$ws.Range("C4").Value = "Lrp1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1.0
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.007329
$ws.Range("H4").Value = 0.021987
$ws.Range("I4").Value = 0.08410957541630165
$ws.Range("J4").Value = 0.08410957541630165
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 249.1329143333333
$ws.Range("N4").Value = 747.398743
$ws.Range("O4").Value = 0.3875740210972192
$ws.Range("P4").Value = 0.3875740210972192
$ws.Range("Q4").Value = 1.825895129149
$ws.Range("R4").Value = 16.433056162341
$ws.Range("S4").Value = 0.03259868635687584
$ws.Range("T4").Value = 0.03259868635687585

# Row 5: ECs -> A2m/Lrp1 -> M2
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "A2m"
$ws.Range("C5").Value = "Lrp1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1.0
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.007329
$ws.Range("H5").Value = 0.021987
$ws.Range("I5").Value = 0.08410957541630165
$ws.Range("J5").Value = 0.08410957541630165
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 186.9310863333334
$ws.Range("N5").Value = 560.793259
$ws.Range("O5").Value = 0.2908071500393791
$ws.Range("P5").Value = 0.2908071500393791
$ws.Range("Q5").Value = 1.370017931737
$ws.Range("R5").Value = 12.330161385633
$ws.Range("S5").Value = 0.0244596659178369
$ws.Range("T5").Value = 0.0244596659178369

# Row 6: ECs -> A2m/Lrp1 -> sCs
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "A2m"
$ws.Range("C6").Value = "Lrp1"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 1.0
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.007329
$ws.Range("H6").Value = 0.021987
$ws.Range("I6").Value = 0.08410957541630165
$ws.Range("J6").Value = 0.08410957541630165
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 17.46044133333334
$ws.Range("N6").Value = 52.38132400000001
$ws.Range("O6").Value = 0.02716306464684043
$ws.Range("P6").Value = 0.02716306464684043
$ws.Range("Q6").Value = 0.127967574532
$ws.Range("R6").Value = 1.151708170788
$ws.Range("S6").Value = 0.002284673834451302
$ws.Range("T6").Value = 0.002284673834451302

# Row 7: FAPs -> A2m/Lrp1 -> ECs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "A2m"
$ws.Range("C7").Value = "Lrp1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 0.07980733333333333
$ws.Range("H7").Value = 0.239422
$ws.Range("I7").Value = 0.9158904245836983
$ws.Range("J7").Value = 0.9158904245836983
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 6.744736333333333
$ws.Range("N7").Value = 20.234209
$ws.Range("O7").Value = 0.01049273071342527
$ws.Range("P7").Value = 0.01049273071342527
$ws.Range("Q7").Value = 0.5382794207997778
$ws.Range("R7").Value = 4.844514787198
$ws.Range("S7").Value = 0.009610191588161484
$ws.Range("T7").Value = 0.009610191588161485

# Row 8: FAPs -> A2m/Lrp1 -> FAPs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "A2m"
$ws.Range("C8").Value = "Lrp1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 0.07980733333333333
$ws.Range("H8").Value = 0.239422
$ws.Range("I8").Value = 0.9158904245836983
$ws.Range("J8").Value = 0.9158904245836983
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 182.5316823333334
$ws.Range("N8").Value = 547.595047
$ws.Range("O8").Value = 0.283963033503136
$ws.Range("P8").Value = 0.2839630335031361
$ws.Range("Q8").Value = 14.56736681587044
$ws.Range("R8").Value = 131.106301342834
$ws.Range("S8").Value = 0.2600790233212622
$ws.Range("T8").Value = 0.2600790233212623

# Row 9: FAPs -> A2m/Lrp1 -> M1
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "A2m"
$ws.Range("C9").Value = "Lrp1"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 0.07980733333333333
$ws.Range("H9").Value = 0.239422
$ws.Range("I9").Value = 0.9158904245836983
$ws.Range("J9").Value = 0.9158904245836983
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 249.1329143333333
$ws.Range("N9").Value = 747.398743
$ws.Range("O9").Value = 0.3875740210972192
$ws.Range("P9").Value = 0.3875740210972192
$ws.Range("Q9").Value = 19.88263353850511
$ws.Range("R9").Value = 178.943701846546
$ws.Range("S9").Value = 0.3549753347403433
$ws.Range("T9").Value = 0.3549753347403434

# Row 10: FAPs -> A2m/Lrp1 -> M2
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "A2m"
$ws.Range("C10").Value = "Lrp1"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 0.07980733333333333
$ws.Range("H10").Value = 0.239422
$ws.Range("I10").Value = 0.9158904245836983
$ws.Range("J10").Value = 0.9158904245836983
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 186.9310863333334
$ws.Range("N10").Value = 560.793259
$ws.Range("O10").Value = 0.2908071500393791
$ws.Range("P10").Value = 0.2908071500393791
$ws.Range("Q10").Value = 14.91847151736645
$ws.Range("R10").Value = 134.266243656298
$ws.Range("S10").Value = 0.2663474841215421
$ws.Range("T10").Value = 0.2663474841215421

# Row 11: FAPs -> A2m/Lrp1 -> sCs
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "A2m"
$ws.Range("C11").Value = "Lrp1"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3.0
$ws.Range("F11").Value = 1.0
$ws.Range("G11").Value = 0.07980733333333333
$ws.Range("H11").Value = 0.239422
$ws.Range("I11").Value = 0.9158904245836983
$ws.Range("J11").Value = 0.9158904245836983
$ws.Range("K11").Value = 3.0
$ws.Range("L11").Value = 1.0
$ws.Range("M11").Value = 17.46044133333334
$ws.Range("N11").Value = 52.38132400000001
$ws.Range("O11").Value = 0.02716306464684043
$ws.Range("P11").Value = 0.02716306464684043
$ws.Range("Q11").Value = 1.393471261636444
$ws.Range("R11").Value = 12.541241354728
$ws.Range("S11").Value = 0.02487839081238912
$ws.Range("T11").Value = 0.02487839081238913
